$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has weekly rating columns B:E (Jun_17, Jun_15, Jun_13, Jun_10),
# newest first. Two newer weeks (Jun_27, Jun_26) are being added, plus an extra
# duplicated "Jun_26" column that mirrors the source edit, so 3 new columns are
# inserted before column B, shifting the existing B:E data (and its cell styles,
# including the highlighted "upgrade" cell) right to E:H.
$ws.Range("B1:D27").Insert(-4161)

# New header labels for the inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns with the default "UN" (unchanged) rating marker for every
# existing analyst row.
$ws.Range("B2:D27").Value = "UN"

# New analyst firms added to the watch list, each starting out at "UN" for the
# three most-recent weeks tracked (B:D); they have no history for the older
# E:H columns.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
